$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $newValue
    $rng.ClearFormats()
}

# Row 2
Set-TextValue "D2" "92.091.56"
Set-TextValue "E2" "  -3.13%  "

# Row 3
Set-TextValue "D3" "3.294.58"
Set-TextValue "E3" "  -5.20%  "

# Row 4
Set-TextValue "E4" "  +0.06%  "

# Row 5
Set-TextValue "D5" "227.41"
Set-TextValue "E5" "  -5.54%  "

# Row 6
Set-TextValue "D6" "609.02"
Set-TextValue "E6" "  -5.72%  "

# Row 7
Set-TextValue "D7" "1.36"
Set-TextValue "E7" "  -7.74%  "

# Row 8
Set-TextValue "D8" "0.378"
Set-TextValue "E8" "  -6.80%  "

# Row 9
Set-TextValue "D9" "0.999"
Set-TextValue "E9" "  -0.03%  "

# Row 10
Set-TextValue "D10" "0.928"
Set-TextValue "E10" "  -8.35%  "

# Row 11
Set-TextValue "D11" "3.295.39"
Set-TextValue "E11" "  -5.16%  "

# Row 12
Set-TextValue "D12" "41.42"
Set-TextValue "E12" "  -2.04%  "

# Row 13
Set-TextValue "E13" "  -3.72%  "

# Row 14
Set-TextValue "D14" "5.92"
Set-TextValue "E14" "  -3.92%  "

# Row 15
Set-TextValue "D15" "91.904.85"
Set-TextValue "E15" "  -3.01%  "

# Row 16
Set-TextValue "D16" "3.902.06"
Set-TextValue "E16" "  -5.18%  "

# Row 17
Set-TextValue "D17" "0.0000240"
Set-TextValue "E17" "  -6.63%  "

# Row 18
Set-TextValue "D18" "7.99"
Set-TextValue "E18" "  -6.41%  "

# Row 19
Set-TextValue "D19" "3.290.90"
Set-TextValue "E19" "  -5.12%  "

# Row 20
Set-TextValue "D20" "17.01"
Set-TextValue "E20" "  -5.55%  "

# Row 21
Set-TextValue "D21" "10.61"
Set-TextValue "E21" "  -7.30%  "

# Row 22
Set-TextValue "E22" "  +5.65%  "

# Row 23
Set-TextValue "D23" "484.36"
Set-TextValue "E23" "  -3.73%  "

# Row 24
Set-TextValue "D24" "0.436"
Set-TextValue "E24" "  -15.59%  "

# Row 25
Set-TextValue "D25" "0.0000177"
Set-TextValue "E25" "  -8.73%  "

# Row 26
Set-TextValue "D26" "6.01"
Set-TextValue "E26" "  -7.74%  "

# Row 27
Set-TextValue "D27" "88.59"
Set-TextValue "E27" "  -3.75%  "

# Row 28
Set-TextValue "D28" "11.61"
Set-TextValue "E28" "  -4.91%  "

# Row 29
Set-TextValue "D29" "3.460.23"
Set-TextValue "E29" "  -5.17%  "

# Row 30
Set-TextValue "E30" "  +0.02%  "

# Row 31
Set-TextValue "E31" "  -7.74%  "

# Row 32
Set-TextValue "E32" "  -2.18%  "

# Row 33
Set-TextValue "D33" "2.58"
Set-TextValue "E33" "  -6.85%  "

# Row 34
Set-TextValue "E34" "  -0.21%  "

# Row 35
Set-TextValue "D35" "0.170"
Set-TextValue "E35" "  -8.01%  "

# Row 36
Set-TextValue "D36" "27.75"
Set-TextValue "E36" "  -10.48%  "

# Row 37
Set-TextValue "E37" "  -9.91%  "

# Row 38
Set-TextValue "D38" "537.94"
Set-TextValue "E38" "  +0.62%  "

# Row 39
Set-TextValue "E39" "  +0.00%  "

# Row 40
Set-TextValue "D40" "7.23"
Set-TextValue "E40" "  -7.66%  "

# Row 41
Set-TextValue "D41" "0.146"
Set-TextValue "E41" "  -3.33%  "

# Row 42
Set-TextValue "E42" "  -7.54%  "

# Row 43
Set-TextValue "D43" "0.850"
Set-TextValue "E43" "  -8.76%  "

# Row 44
Set-TextValue "D44" "23.75"
Set-TextValue "E44" "  -1.41%  "

# Row 45
Set-TextValue "B45" "ImmutableX"
Set-TextValue "C45" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D45" "1.65"
Set-TextValue "E45" "  -3.20%  "

# Row 46
Set-TextValue "B46" "MantraDAO"
Set-TextValue "C46" "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
Set-TextValue "D46" "3.57"
Set-TextValue "E46" "  +1.39%  "

# Row 47
Set-TextValue "E47" "  -4.19%  "

# Row 48
Set-TextValue "D48" "5.25"
Set-TextValue "E48" "  -8.31%  "

# Row 49
Set-TextValue "D49" "2.06"
Set-TextValue "E49" "  -4.82%  "

# Row 50
Set-TextValue "D50" "51.42"
Set-TextValue "E50" "  -3.64%  "

# Row 51
Set-TextValue "D51" "7.83"
Set-TextValue "E51" "  -3.38%  "
